$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colC = @(0.1797227757417943,0.1754687248507025,0.1729329611819423,0.1719188298279164,0.1717515957463718,0.1729192064027956,0.1782401782457157,0.1892784446874174,0.1977560415493258,0.2016925905216738,0.2031947485272951,0.2028707222923174,0.2018159442599199,0.2011713551705867,0.1975003862764879,0.1952688403000735,0.1939928507872963,0.1935621182829408,0.1955056124378416,0.20212544708167,0.2065187440495606,0.2041678542607315,0.1953985460864658,0.1862276927170967)
$arrC = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrC[$i,0] = $colC[$i] }
$ws.Range("C2:C25").Value = $arrC

$colD = @(0.0435955048535881,0.04360659440974857,0.04361990373555358,0.04362696990906656,0.04362824268430998,0.04361999237159608,0.043597984053104,0.04360608034268765,0.04364288194835453,0.04366624766657878,0.04367604204119147,0.04367389065857097,0.04366703452650356,0.04366295799954756,0.04364148769952791,0.04363000846613119,0.04362402980011382,0.04362211287618933,0.04363116592786298,0.04366902270509243,0.04369927725776535,0.04368262740190687,0.04363064070457412,0.04359844471111529)
$arrD = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrD[$i,0] = $colD[$i] }
$ws.Range("D2:D25").Value = $arrD

$colE = @(0.1325008509881194,0.1292269150301415,0.1272877263987304,0.1265153371281009,0.1263881599602072,0.1272772374459343,0.1313572480471592,0.1399228524862224,0.1465627963920042,0.1496593895848974,0.1508429568416432,0.1505875669401888,0.1497565427184071,0.149248943571429,0.1463619586457767,0.144610378157374,0.1436100723548677,0.1432726146511243,0.144796096036977,0.1500003372183443,0.1534654707811001,0.1516102141028597,0.1447121121127211,0.1375449137789175)
$arrE = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrE[$i,0] = $colE[$i] }
$ws.Range("E2:E25").Value = $arrE

$colF = @(1.204450486139535,1.195608402036711,1.191054737370393,1.189418574994605,1.189160129785805,1.19103178364098,1.201219653845612,1.228175713859628,1.252282289933135,1.264193836553318,1.268841161941708,1.267834186466047,1.264573431248309,1.262593946554531,1.251522928390898,1.244973868907977,1.241295958750214,1.240065933816851,1.245661817632424,1.265527478671316,1.279307867106468,1.271879839668273,1.245350524353626,1.220131415021442)
$arrF = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrF[$i,0] = $colF[$i] }
$ws.Range("F2:F25").Value = $arrF

$colG = @(1.078133462989925,1.070244041303539,1.066257283302733,1.064847430968072,1.064626274834623,1.066237400966628,1.075234746704368,1.099720704628723,1.121939975585889,1.132979332253342,1.137294602034927,1.136359216707774,1.133331643064622,1.131494763030645,1.1212373550199,1.115184056338222,1.111789999174079,1.110655853821186,1.115819363210875,1.134217246128543,1.14702809020045,1.140118404682767,1.115531872961668,1.092358098250585)
$arrG = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrG[$i,0] = $colG[$i] }
$ws.Range("G2:G25").Value = $arrG

$colH = @(1.050730821794829,1.053200183601831,1.055334052207968,1.056358695817693,1.056538196180028,1.055347243346247,1.051453910802692,1.048732201971163,1.049746282722367,1.050866180750546,1.05138527846097,1.051269250574677,1.050906981154668,1.050697463979816,1.049686373854627,1.049235004147192,1.049037358729407,1.048981073406111,1.049276637279291,1.051010807231677,1.052698164763882,1.05174679716589,1.049257622295499,1.048940608978825)
$arrH = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrH[$i,0] = $colH[$i] }
$ws.Range("H2:H25").Value = $arrH

$colI = @(1.183615314825666,1.177843412716641,1.175109838312117,1.174199151557183,1.174060194592499,1.175096734136353,1.181456649605714,1.200384275800431,1.218265327799813,1.227271753363482,1.230808329492945,1.230041049187179,1.2275601801025,1.226057008252994,1.217694340581787,1.212787950543188,1.210047987848498,1.20913435944513,1.213301746534441,1.228285445958704,1.238813093824419,1.233126840534979,1.213069207772065,1.194568667089527)
$arrI = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrI[$i,0] = $colI[$i] }
$ws.Range("I2:I25").Value = $arrI

$colK = @(1.662797164383733,1.510329207646521,1.417070023065776,1.379156807389847,1.372866850861215,1.416558343833515,1.610152684091133,1.992600196344597,2.275299330270911,2.404281518418713,2.453178061720848,2.442644939572574,2.408303194076154,2.387274859085665,2.266877621246294,2.193114671466901,2.15072408700911,2.136377562106873,2.200963153730584,2.41838874518993,2.560802009675285,2.484765046222265,2.197414802752121,1.888837073798356)
$arrK = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrK[$i,0] = $colK[$i] }
$ws.Range("K2:K25").Value = $arrK

$colL = @(0.1772876741929679,0.1727706243253593,0.1701043161996765,0.1690446634995126,0.1688703312619424,0.1700899165578207,0.1757079113993285,0.1875787959966431,0.1968271391290415,0.2011502561716867,0.2028040801581454,0.2024471532908194,0.2012859810828758,0.200576913541056,0.1965469536694684,0.1941044655401782,0.192710524967481,0.1922404329266811,0.1943633424887423,0.2016265902404228,0.2064712285825152,0.2038765913195562,0.1942462722115152,0.1842753443662275)
$arrL = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrL[$i,0] = $colL[$i] }
$ws.Range("L2:L25").Value = $arrL

$colN = @(1.236459234024103,1.248587195801505,1.256560830671638,1.259942714423005,1.260512281433996,1.256605903145264,1.240531595194184,1.213190026043563,1.195648712027293,1.188221475815354,1.18548842081367,1.186073497498846,1.187995032591282,1.189182380180384,1.196145218589677,1.200558205037979,1.203148429033099,1.204034362639447,1.200083053959062,1.18742847374596,1.179621260121053,1.183745709431349,1.200297704239325,1.220139351523045)
$arrN = New-Object 'object[,]' 24,1
for ($i = 0; $i -lt 24; $i++) { $arrN[$i,0] = $colN[$i] }
$ws.Range("N2:N25").Value = $arrN
